# Remove the two "MARCO" balance rows (accounts 004435987 and 004436055)
# from the "Export" sheet. Rows are deleted from the bottom up so that the
# row index of the earlier row is not shifted by the later deletion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: 004436055 / MARCO / 19911.52
$ws.Rows(5).Delete()

# Row 3: 004435987 / MARCO / 32941.09
$ws.Rows(3).Delete()
